$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '43.561.91'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.88%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.386.18'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +5.94%  '
$ws.Range('E4').Value = '  -0.11%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '236.87'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.97%  '
$ws.Range('E6').Value = '  +2.14%  '
$ws.Range('E7').Value = '  +13.58%  '
$ws.Range('E8').Value = '  +0.01%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.468'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +4.06%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0976'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.40%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '57.33'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -3.57%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '26.75'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +0.82%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '2.727.97'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +5.48%  '
$ws.Range('E14').Value = '  -0.08%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '15.88'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +1.94%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '6.28'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +2.92%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.860'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +3.07%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.375.95'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +5.52%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '43.441.94'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.88%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.0₃0993'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.65%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.43'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +5.57%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '74.73'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.89%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '251.75'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +2.02%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '3.92'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +17.49%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.04%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.50'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +3.16%  '
$ws.Range('E27').Value = '  +2.64%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '23.04'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +6.17%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '10.05'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.83%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '174.55'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.68%  '
$ws.Range('E31').Value = '  +7.86%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.128'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -7.82%  '
$ws.Range('E33').Value = '  +1.18%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.03'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +4.16%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.0693'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.95%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '5.10'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +3.21%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.47'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +7.77%  '
$ws.Range('E38').Value = '  +2.86%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '3.70'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('E40').Value = '  +0.89%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '18.99'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +10.77%  '
$ws.Range('B42').Value = 'BinanceUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '8.91'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +3.56%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.19'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +10.37%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '4.57'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +4.01%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '100.08'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +2.06%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.23'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.41%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0953'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.46%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.454.47'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.78%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.599.62'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +5.63%  '
$ws.Range('E51').Value = '  -0.37%  '
